$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.006.74'
$ws.Range('E2').Value = '  +2.59%  '
$ws.Range('D3').Value = '2.625.02'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'601.13"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.29%  '
$ws.Range('D6').Value = "'154.78"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.75%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = "'0.586"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').Value = "'0.118"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +8.65%  '
$ws.Range('D10').Value = "'0.408"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +5.88%  '
$ws.Range('D11').Value = "'5.75"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').Value = "'0.154"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.74%  '
$ws.Range('D13').Value = "'29.16"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +6.09%  '
$ws.Range('D14').Value = "'0.0000189"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +22.35%  '
$ws.Range('D15').Value = '3.094.34'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '64.901.90'
$ws.Range('D17').Value = '2.609.10'
$ws.Range('D18').Value = "'12.53"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +4.38%  '
$ws.Range('D19').Value = "'4.90"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.42%  '
$ws.Range('D20').Value = "'358.20"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.36%  '
$ws.Range('E21').Value = '  +7.58%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = "'68.68"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.35%  '
$ws.Range('D24').Value = "'1.65"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.89%  '
$ws.Range('D25').Value = "'9.37"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.17%  '
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('E27').Value = '  +2.74%  '
$ws.Range('D28').Value = "'8.15"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '0.0₃0945'
$ws.Range('E30').Value = '  +12.75%  '
$ws.Range('E31').Value = '  +5.26%  '
$ws.Range('D32').Value = "'523.58"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -7.49%  '
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('D34').Value = "'5.47"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.53%  '
$ws.Range('D35').Value = "'6.33"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +5.55%  '
$ws.Range('D36').Value = "'0.425"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('E37').Value = '  +5.26%  '
$ws.Range('D38').Value = "'162.15"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('D39').Value = "'2.01"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.30%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = "'42.15"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.69%  '
$ws.Range('D43').Value = "'164.34"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').Value = "'4.14"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.41%  '
$ws.Range('E45').Value = '  +5.82%  '
$ws.Range('D46').Value = "'23.14"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.35%  '
$ws.Range('D47').Value = "'2.21"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +5.51%  '
$ws.Range('D48').Value = "'0.652"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.58%  '
$ws.Range('D49').Value = "'0.0263"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +6.89%  '
$ws.Range('D50').Value = "'0.0979"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('D51').Value = "'19.47"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.46%  '
